$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.338.47"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.932.83"
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  +0.23%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "251.22"
$ws.Range("E5").Value = "  +1.96%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.7178"
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +0.24%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3278"
$ws.Range("E8").Value = "  +1.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "27.56"
$ws.Range("E9").Value = "  +4.92%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07191"
$ws.Range("E10").Value = "  +5.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8026"
$ws.Range("E11").Value = "  +1.06%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08079"
$ws.Range("E12").Value = "  +1.99%  "
$ws.Range("D13").Value = "1.931.93"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.418"
$ws.Range("E14").Value = "  +0.62%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "94.58"
$ws.Range("E15").Value = "  +0.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.90"
$ws.Range("E16").Value = "  +3.24%  "
$ws.Range("D17").Value = "30.333.36"
$ws.Range("E17").Value = "  +0.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "252.58"
$ws.Range("E18").Value = "  -2.42%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008141"
$ws.Range("E19").Value = "  +2.96%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.806"
$ws.Range("E20").Value = "  -0.33%  "
$ws.Range("D21").Value = "2.187.33"
$ws.Range("E21").Value = "  +0.81%  "
$ws.Range("E22").Value = "  +0.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.936"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.729"
$ws.Range("E25").Value = "  +0.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.30"
$ws.Range("E26").Value = "  +3.01%  "
$ws.Range("B27").Value = "EthereumClassic"
$ws.Range("C27").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.24"
$ws.Range("E27").Value = "  +2.54%  "
$ws.Range("B28").Value = "LidoDAOToken"
$ws.Range("C28").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.334"
$ws.Range("E28").Value = "  +4.44%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1292"
$ws.Range("E29").Value = "  -3.19%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.360"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.544"
$ws.Range("E31").Value = "  -0.02%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.425"
$ws.Range("E32").Value = "  +0.38%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.207"
$ws.Range("E33").Value = "  +0.59%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05210"
$ws.Range("E34").Value = "  +3.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.264"
$ws.Range("E35").Value = "  +6.15%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7473"
$ws.Range("E36").Value = "  +1.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.764"
$ws.Range("E37").Value = "  +1.22%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01964"
$ws.Range("E38").Value = "  +1.03%  "
$ws.Range("E39").Value = "  -0.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.98"
$ws.Range("E40").Value = "  -1.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.470"
$ws.Range("E41").Value = "  -0.16%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4533"
$ws.Range("E42").Value = "  +2.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.026"
$ws.Range("E43").Value = "  +1.05%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8404"
$ws.Range("E45").Value = "  +1.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.97"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.789"
$ws.Range("E47").Value = "  +1.10%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.413"
$ws.Range("E48").Value = "  +2.17%  "
$ws.Range("E49").Value = "  +2.01%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.4180"
$ws.Range("E50").Value = "  +2.06%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06062"
$ws.Range("E51").Value = "  +2.54%  "
